# Update footprint + schema to use cheaper CP2102
# Rewrites the full component-placement table (rows 2-47) on the active sheet
# to the final designator / Mid X / Mid Y / Rotation layout. The row count
# grows from 39 to 46 data rows (dimension A1:E40 -> A1:E47).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the numeric-cell formatting (style "0.000000") used by the existing
# B:D data rows down to the newly added rows before filling in values.
$ws.Range("B2:D2").Copy() | Out-Null
$ws.Range("B41:D47").PasteSpecial(-4122) | Out-Null

$data = @(
    @("C1", 142.6, -85.2, 90.0),
    @("C2", 134.5125, -83.05, 0.0),
    @("C3", 151.55, -85.0, 90.0),
    @("C4", 154.75, -85.65, 180.0),
    @("C5", 149.1, -90.8, 180.0),
    @("C6", 156.1, -69.9625, 90.0),
    @("C7", 151.1, -96.0, 180.0),
    @("C8", 149.3, -86.1, 90.0),
    @("C9", 158.4, -71.0, -90.0),
    @("C10", 152.7, -88.5, 180.0),
    @("C11", 156.5875, -88.5, 0.0),
    @("C12", 152.625, -99.05, 0.0),
    @("C13", 117.85, -116.9, -90.0),
    @("C14", 131.4, -120.4, -90.0),
    @("C15", 117.85, -110.2, 180.0),
    @("C16", 157.1625, -123.27, 0.0),
    @("C17", 129.6, -78.8, 180.0),
    @("C18", 134.7, -78.8, 0.0),
    @("D1", 154.15, -90.85, 180.0),
    @("D2", 156.55, -93.95, 180.0),
    @("D3", 151.65, -93.45, 0.0),
    @("D4", 137.471501, -72.37, 90.0),
    @("J1", 146.35, -77.7, 180.0),
    @("J2", 124.8, -120.75, 0.0),
    @("J4", 119.4, -89.0, 180.0),
    @("JP1", 130.35, -94.53, 0.0),
    @("L1", 156.65, -98.05, -90.0),
    @("Q1", 138.65, -86.5, 0.0),
    @("Q2", 154.65, -119.85, 180.0),
    @("Q3", 149.3, -102.24, 180.0),
    @("R1", 134.9, -85.85, 180.0),
    @("R2", 158.155, -120.82, 180.0),
    @("R3", 135.3, -89.53, 0.0),
    @("R4", 158.91, -118.49, -90.0),
    @("R5", 152.92, -103.21, 180.0),
    @("R6", 117.65, -120.75, -90.0),
    @("R7", 127.5, -114.1, 0.0),
    @("R8", 155.0, -115.85, 90.0),
    @("R9", 121.025, -95.82, 0.0),
    @("R10", 118.825, -94.2, 180.0),
    @("R11", 128.95, -90.55, 0.0),
    @("R12", 123.15, -88.35, -90.0),
    @("U1", 142.514, -114.75, 180.0),
    @("U2", 124.9, -108.1, 90.0),
    @("U3", 121.7, -113.0625, -90.0),
    @("U4", 127.85, -85.5, 0.0)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = "top"
    $row++
}
